$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 333, shifting existing rows 333:400 down to 334:401
# (this also extends the sheet dimension from A1:T400 to A1:T401 and
# keeps every shifted row's formatting, e.g. the date style on column D)
$ws.Rows(333).Insert()

# Populate the newly inserted row 333 with the new day's record.
# Columns A,B,C,E,F,G,H,I,J,K,T are constant across this block of rows.
$ws.Range("A333").Value = 5
$ws.Range("B333").Value = "Macroferia Regional de Talca"
$ws.Range("C333").Value = "Maule"
$ws.Range("D333").Value = 44543
$ws.Range("E333").Value = 7
$ws.Range("F333").Value = "Fruta"
$ws.Range("G333").Value = 100101
$ws.Range("H333").Value = "Berries"
$ws.Range("I333").Value = 100112025
$ws.Range("J333").Value = "Frutilla"
$ws.Range("K333").Value = "Sin especificar"
$ws.Range("L333").Value = "Especial"
$ws.Range("M333").Value = 450
$ws.Range("N333").Value = 7500
$ws.Range("O333").Value = 8000
$ws.Range("P333").Value = 7722
$ws.Range("Q333").Value = '$/bandeja 7 kilos'
$ws.Range("R333").Value = "Provincia de Melipilla"
$ws.Range("S333").Value = 1103
$ws.Range("T333").Value = 7
